$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HoaDonMau")

# Helper: write a numeric-looking string into a cell as literal text while
# preserving that cell's pre-existing style (a plain ".Value = '60,600'"
# assignment would get auto-coerced to a real number by Excel). We stage the
# text in an unused scratch cell formatted as Text, copy it, and paste
# values-only into the destination so the destination keeps its own
# formatting/style untouched.
$scratch = $ws.Range("Z1")
function Set-LiteralText($range, $text) {
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)
    $scratch.Clear()
}

# Header info: invoice number, date/time, salesperson
$ws.Range("A6").Value = "Số HD: HD131223005"
$ws.Range("A7").Value = "Ngày giờ: 06:26:37 - 13/12/2023"
$ws.Range("A9").Value = "NV bán hàng: Dương Thái Bảo"

# Line item: product name, unit price, quantity, line total
$ws.Range("B12").Value = "Mắt biếc"
Set-LiteralText $ws.Range("B13") "60,600"
$ws.Range("C13").Value = 7
Set-LiteralText $ws.Range("D13") "424,200"

# Totals block
Set-LiteralText $ws.Range("D14") "424,200"
Set-LiteralText $ws.Range("D15") "84,840"
Set-LiteralText $ws.Range("D16") "339,360"
Set-LiteralText $ws.Range("D17") "339,360"
